# Adds new team tasks (rows 13-17) to the Tasks sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 : Task 12 -------------------------------------------------
# (Values are written in a particular column order so that new shared
#  strings land in the same order as in the target workbook.)
$ws.Range("A13").Value = 12
$ws.Range("D13").Value = "Task 5 is finished"
$ws.Range("F13").Value = "New"
$ws.Range("B13").Value = "Filled historicaldata (solutions and tickets) in SupportTicketingSystem/BackendSystem/HISTORYDATA.xlsx"
$ws.Range("C13").Value = "Marcia"
$ws.Range("E13").Value = "1. team review 2. uploaded to github 3. present to a team"

# --- Row 14 : Task 13 --------------------------------------------------
$ws.Range("A14").Value = 13
$ws.Range("E14").Value = "1. team review"
$ws.Range("D14").Value = "Tasks 1-11 are finished"

# Row 17's "Prerequisites" text is entered now so the shared-string table
# keeps the same ordering as the target workbook.
$ws.Range("D17").Value = "Tasks 1-12 are finished"

$ws.Range("B14").Value = "Check ODATA API with POSTMAN. Do we have all ODATA functionality regarding our mockup?"
$ws.Range("C14").Value = "Shamil"
$ws.Range("F14").Value = "New"

# --- Row 15 : Task 14 ---------------------------------------------------
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Check ODATA API with POSTMAN. Do we have all ODATA functionality regarding our mockup?"
$ws.Range("C15").Value = "Ilkay"
$ws.Range("D15").Value = "Tasks 1-11 are finished"
$ws.Range("E15").Value = "1. team review"
$ws.Range("F15").Value = "New"

# --- Row 16 : Task 15 ---------------------------------------------------
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Check ODATA API with POSTMAN. Do we have all ODATA functionality regarding our mockup?"
$ws.Range("C16").Value = "Marcus"
$ws.Range("D16").Value = "Tasks 1-11 are finished"
$ws.Range("E16").Value = "1. team review"
$ws.Range("F16").Value = "New"

# --- Row 17 : Task 16 (remaining cells) ---------------------------------
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Check ODATA API with POSTMAN. Do we have all ODATA functionality regarding our mockup?"
$ws.Range("C17").Value = "Marcia"
$ws.Range("E17").Value = "1. team review"
$ws.Range("F17").Value = "New"

# Match the wrap-text formatting (columns A:E) used by every other data row.
$ws.Range("A13:E17").WrapText = $true

# Rows hold two-line wrapped text, same as other long-text rows in the sheet.
$ws.Range("13:17").RowHeight = 28.8

# Put the selection where the author left it after typing the new tasks.
$ws.Range("B14").Select()
